$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = 0
$ws.Range("E6").Value = 442.577832503079
$ws.Range("H6").Value = 0
$ws.Range("K6").Value = 432.18685872638656
$ws.Range("N6").Value = 427.0098623277617
$ws.Range("Q6").Value = 395.8726366256971
$ws.Range("T6").Value = 264.6000743293567
$ws.Range("W6").Value = 415.44214382305375

$ws.Range("B7").Value = 0
$ws.Range("E7").Value = 1.2470766179992978
$ws.Range("H7").Value = 0
$ws.Range("K7").Value = 1.2379625509763041
$ws.Range("N7").Value = 1.2302053288931263
$ws.Range("Q7").Value = 1.210508482108495
$ws.Range("T7").Value = 1.2707780553297445
$ws.Range("W7").Value = 1.2126185687161064

$ws.Range("B8").Value = 0
$ws.Range("E8").Value = 131.02929347196422
$ws.Range("H8").Value = 0
$ws.Range("K8").Value = 169.9948920291913
$ws.Range("N8").Value = 191.06040697591587
$ws.Range("Q8").Value = 270.11168062166183
$ws.Range("T8").Value = 684.1505966703487
$ws.Range("W8").Value = 219.14486799236963

$ws.Range("B9").Value = 0
$ws.Range("E9").Value = -0.2411028922254593
$ws.Range("H9").Value = 0
$ws.Range("K9").Value = -2.585499296861247
$ws.Range("N9").Value = -0.13734193778968887
$ws.Range("Q9").Value = -0.0937559271795416
$ws.Range("T9").Value = -0.13748416564420185
$ws.Range("W9").Value = 0.03184726528250274

$ws.Range("B10").Value = 0
$ws.Range("E10").Value = 1.0257594035284696
$ws.Range("H10").Value = 0
$ws.Range("K10").Value = 1.1358487001869144
$ws.Range("N10").Value = 1.0079394810821154
$ws.Range("Q10").Value = 0.959131926706989
$ws.Range("T10").Value = 1.1611485907425028
$ws.Range("W10").Value = 0.9277455470066464

$ws.Range("B11").Value = 0
$ws.Range("E11").Value = 0.0005229739435035055
$ws.Range("H11").Value = 0
$ws.Range("K11").Value = 0.0058407473290235925
$ws.Range("N11").Value = 0.0003129402848704141
$ws.Range("Q11").Value = 0.00023704154085226256
$ws.Range("T11").Value = 0.00034410740708201937
$ws.Range("W11").Value = -0.00007048825820544594

$ws.Range("B12").Value = 0
$ws.Range("E12").Value = -0.0001277359011035134
$ws.Range("H12").Value = 0
$ws.Range("K12").Value = -0.0002879212469772541
$ws.Range("N12").Value = -0.00006095177807661652
$ws.Range("Q12").Value = -0.0000343168459982109
$ws.Range("T12").Value = -0.00025868165908548386
$ws.Range("W12").Value = -0.00002656007442604687

$ws.Range("B13").Value = 0
$ws.Range("E13").Value = -0.0000002817421002053448
$ws.Range("H13").Value = 0
$ws.Range("K13").Value = -0.0000032981190174523347
$ws.Range("N13").Value = -0.00000017791791555844049
$ws.Range("Q13").Value = -0.00000015018590181111645
$ws.Range("T13").Value = -0.0000002121354047318504
$ws.Range("W13").Value = 0.000000038277878666607215

$ws.Range("B14").Value = 0
$ws.Range("E14").Value = 1.0028562195969548
$ws.Range("H14").Value = 0
$ws.Range("K14").Value = 1.0486056180704453
$ws.Range("N14").Value = 1.0718793376412286
$ws.Range("Q14").Value = 1.0432674461991223
$ws.Range("T14").Value = 2.271740408941285
$ws.Range("W14").Value = 1.0493812633124497

$ws.Range("B15").Value = 0
$ws.Range("E15").Value = -0.00001205476262848713
$ws.Range("H15").Value = 0
$ws.Range("K15").Value = -0.00011236041808676488
$ws.Range("N15").Value = -0.00017184062858964448
$ws.Range("Q15").Value = -0.00010496184607741685
$ws.Range("T15").Value = -0.0027848983050172724
$ws.Range("W15").Value = -0.00012197170060693368

$ws.Range("B16").Value = 0
$ws.Range("E16").Value = 0.22910582470681978
$ws.Range("H16").Value = 0
$ws.Range("K16").Value = 0.28941234581916875
$ws.Range("N16").Value = -0.03087456709901914
$ws.Range("Q16").Value = 0.17944568252976265
$ws.Range("T16").Value = 0.06550371643571735
$ws.Range("W16").Value = 0.21493003381173162

$ws.Range("B17").Value = 0
$ws.Range("E17").Value = 0.000000009687514556643437
$ws.Range("H17").Value = 0
$ws.Range("K17").Value = 0.00000006457210417156402
$ws.Range("N17").Value = 0.00000010229935790416947
$ws.Range("Q17").Value = 0.00000006363392607075201
$ws.Range("T17").Value = 0.000001524244345134667
$ws.Range("W17").Value = 0.00000007505965809255445

$ws.Range("B18").Value = 0
$ws.Range("E18").Value = -0.0001712805315216371
$ws.Range("H18").Value = 0
$ws.Range("K18").Value = -0.0002099507396073081
$ws.Range("N18").Value = 0.000018868118559906043
$ws.Range("Q18").Value = -0.00014867106419065064
$ws.Range("T18").Value = -0.000075647147805293
$ws.Range("W18").Value = -0.00017219647310054168

$ws.Range("B19").Value = 0
$ws.Range("E19").Value = -0.5345656677137223
$ws.Range("H19").Value = 0
$ws.Range("K19").Value = -0.5042452047504
$ws.Range("N19").Value = -0.45675813176519686
$ws.Range("Q19").Value = -0.3826111956090942
$ws.Range("T19").Value = -0.5041951645729208
$ws.Range("W19").Value = -0.42658800672376473
